$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1764705882352941
$ws.Range("C2").Value = 0.583710407239819
$ws.Range("J2").Value = 0.01357466063348416
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.07239819004524888
$ws.Range("C3").Value = 0.05185185185185185
$ws.Range("J3").Value = 0.02222222222222222
$ws.Range("P3").Value = 0.7555555555555555
$ws.Range("S3").Value = 0.1703703703703704
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.5581395348837209
$ws.Range("S4").Value = 0.3720930232558139
$ws.Range("B6").Value = 0.05853658536585366
$ws.Range("D6").Value = 0.01463414634146342
$ws.Range("F6").Value = 0.1219512195121951
$ws.Range("J6").Value = 0.2
$ws.Range("O6").Value = 0.01951219512195122
$ws.Range("Q6").Value = 0.1658536585365854
$ws.Range("R6").Value = 0.06829268292682927
$ws.Range("S6").Value = 0.3512195121951219
$ws.Range("B7").Value = 0.08762886597938144
$ws.Range("D7").Value = 0.02061855670103093
$ws.Range("E7").Value = 0.005154639175257732
$ws.Range("F7").Value = 0.06185567010309279
$ws.Range("J7").Value = 0.1030927835051546
$ws.Range("O7").Value = 0.05670103092783505
$ws.Range("Q7").Value = 0.2010309278350516
$ws.Range("R7").Value = 0.05670103092783505
$ws.Range("S7").Value = 0.4072164948453608
$ws.Range("B8").Value = 0.07631578947368421
$ws.Range("D8").Value = 0.01842105263157895
$ws.Range("F8").Value = 0.05
$ws.Range("J8").Value = 0.1289473684210526
$ws.Range("O8").Value = 0.02368421052631579
$ws.Range("Q8").Value = 0.1815789473684211
$ws.Range("R8").Value = 0.1
$ws.Range("S8").Value = 0.4210526315789473
$ws.Range("B9").Value = 0.07920792079207921
$ws.Range("D9").Value = 0.0198019801980198
$ws.Range("F9").Value = 0.08415841584158416
$ws.Range("J9").Value = 0.0891089108910891
$ws.Range("O9").Value = 0.0297029702970297
$ws.Range("Q9").Value = 0.1386138613861386
$ws.Range("R9").Value = 0.09405940594059406
$ws.Range("S9").Value = 0.4653465346534654
$ws.Range("B10").Value = 0.1048951048951049
$ws.Range("D10").Value = 0.02497502497502498
$ws.Range("E10").Value = 0.000999000999000999
$ws.Range("F10").Value = 0.07292707292707293
$ws.Range("J10").Value = 0.1358641358641359
$ws.Range("O10").Value = 0.02297702297702298
$ws.Range("Q10").Value = 0.2307692307692308
$ws.Range("R10").Value = 0.06693306693306693
$ws.Range("S10").Value = 0.3396603396603396
$ws.Range("G11").Value = 0.1458333333333333
$ws.Range("J11").Value = 0.04583333333333333
$ws.Range("K11").Value = 0.1625
$ws.Range("L11").Value = 0.6333333333333333
$ws.Range("S11").Value = 0.0125
$ws.Range("G12").Value = 0.8170731707317073
$ws.Range("J12").Value = 0.1158536585365854
$ws.Range("K12").Value = 0.006097560975609756
$ws.Range("L12").Value = 0.05487804878048781
$ws.Range("S12").Value = 0.006097560975609756
$ws.Range("G13").Value = 0.7878787878787878
$ws.Range("J13").Value = 0.1818181818181818
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("F15").Value = 0.02439024390243903
$ws.Range("H15").Value = 0.1707317073170732
$ws.Range("I15").Value = 0.06829268292682927
$ws.Range("J15").Value = 0.3414634146341464
$ws.Range("K15").Value = 0.03902439024390244
$ws.Range("M15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.08780487804878048
$ws.Range("S15").Value = 0.2634146341463415
$ws.Range("F16").Value = 0.01282051282051282
$ws.Range("H16").Value = 0.2115384615384615
$ws.Range("I16").Value = 0.1153846153846154
$ws.Range("J16").Value = 0.3782051282051282
$ws.Range("K16").Value = 0.07692307692307693
$ws.Range("M16").Value = 0.02564102564102564
$ws.Range("O16").Value = 0.0576923076923077
$ws.Range("S16").Value = 0.1217948717948718
$ws.Range("F17").Value = 0.01256281407035176
$ws.Range("H17").Value = 0.1482412060301508
$ws.Range("I17").Value = 0.1155778894472362
$ws.Range("J17").Value = 0.4170854271356784
$ws.Range("K17").Value = 0.1256281407035176
$ws.Range("M17").Value = 0.02010050251256281
$ws.Range("N17").Value = 0.002512562814070352
$ws.Range("O17").Value = 0.06281407035175879
$ws.Range("S17").Value = 0.09547738693467336
$ws.Range("F18").Value = 0.0272108843537415
$ws.Range("H18").Value = 0.2312925170068027
$ws.Range("I18").Value = 0.08163265306122448
$ws.Range("J18").Value = 0.3673469387755102
$ws.Range("K18").Value = 0.08843537414965986
$ws.Range("M18").Value = 0.02040816326530612
$ws.Range("N18").Value = 0.006802721088435374
$ws.Range("O18").Value = 0.07482993197278912
$ws.Range("S18").Value = 0.1020408163265306
$ws.Range("F19").Value = 0.01684836471754212
$ws.Range("H19").Value = 0.2091179385530228
$ws.Range("I19").Value = 0.1080277502477701
$ws.Range("J19").Value = 0.3716551040634292
$ws.Range("K19").Value = 0.1149653121902874
$ws.Range("M19").Value = 0.01684836471754212
$ws.Range("O19").Value = 0.09018830525272548
$ws.Range("S19").Value = 0.07234886025768088
